# "Transform to MVC model" - fill in previously-blank data cells on Sheet1
# (rows 4-7) with the values that belong there, and move the active-cell
# selection from F18 to M18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Drug3): generic/brand name + manufacturer were blank
$ws.Range("B4").Value = "Anastrazole"
$ws.Range("F4").Value = "Anastrazole"
$ws.Range("G4").Value = "Ranbaxy limited (a Sun Pharmaceuticals company)"

# Row 5 (Drug4): strength unit, brand name, AMC were blank
$ws.Range("E5").Value = "mg"
$ws.Range("F5").Value = "Anastrazole"
$ws.Range("L5").Value = 432

# Row 6 (Drug5): available stock was blank
$ws.Range("K6").Value = 45

# Row 7 (Drug6): form, expiry date, AMC, monthly consumption were blank
$ws.Range("C7").Value = "Tab"
$ws.Range("I7").Value = 42715
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 4

# Update the saved selection to M18
$ws.Range("M18").Select()
